# Reworked on the project as said.
$wb = $excel.ActiveWorkbook

# The "php" sheet holds the Departure/Destination values that need trimming.
$ws = $wb.Worksheets.Item("php")
$ws.Activate()

# Remove the leading space from the Departure (B2) and Destination (C2) values.
$ws.Range("B2").Value = "Chennai"
$ws.Range("C2").Value = "McCarran"

# Move the active selection from F1 to B2.
$ws.Range("B2").Select()
